# Daily refresh of the BP Terminal Gate Pricing (TGP) workbook.
#
# Each effective-date block rolls forward by one day (the "yesterday" row
# becomes "today", a new "today" row is appended with refreshed pricing,
# and the D/E/F/G price columns are updated with the latest cents-per-litre
# figures). Values below are the literal post-edit figures taken from the
# daily pricing feed, applied cell-by-cell so existing number formats /
# cell styles (dates in column A, 2-decimal prices in D:G) are preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A, D, E, F, G (in that
# order). A column that doesn't apply to a given row (no data, or an "N/A"
# text cell that isn't touched by this edit) is simply omitted.
$updates = @(
    @{ Row = 8;  A = 46046; D = 158.15;              E = 150.47999999999999; F = 160.47999999999999; G = 150.5 }
    @{ Row = 9;  A = 46046; D = 158.15;              E = 150.47999999999999; F = 160.47999999999999; G = 150.5 }
    @{ Row = 10; A = 46046; D = 158.94;              E = 152.26;             F = 162.26;             G = 152.66 }
    @{ Row = 11; A = 46045; D = 157.83000000000001;  E = 150.65;             F = 160.65;             G = 150.66 }
    @{ Row = 12; A = 46045; D = 157.83000000000001;  E = 150.65;             F = 160.65;             G = 150.66 }
    @{ Row = 13; A = 46045; D = 158.59;              E = 152.4;              F = 162.4;              G = 152.80000000000001 }

    @{ Row = 17; A = 46046; D = 162.34;              E = 154.77000000000001; F = 164.77 }
    @{ Row = 18; A = 46045; D = 161.97;              E = 154.88;             F = 164.88 }

    @{ Row = 22; A = 46046; D = 159.22;              E = 152.05000000000001; F = 161.65;             G = 153.12 }
    @{ Row = 23; A = 46046; D = 163.72;              E = 157.59;             F = 167.59 }
    @{ Row = 24; A = 46046; D = 163.85;              E = 158.31;             F = 168.31 }
    @{ Row = 25; A = 46046; D = 163.84;              E = 157.85;             F = 167.85;             G = 157.97999999999999 }
    @{ Row = 26; A = 46046; D = 163.41999999999999;  E = 159.47;             F = 169.47 }
    @{ Row = 27; A = 46045; D = 158.91;              E = 152.22;             F = 161.82;             G = 153.29 }
    @{ Row = 28; A = 46045; D = 163.37;              E = 157.72999999999999; F = 167.73 }
    @{ Row = 29; A = 46045; D = 163.5;               E = 158.46;             F = 168.46 }
    @{ Row = 30; A = 46045; D = 163.49;              E = 158;                F = 168;                G = 158.13 }
    @{ Row = 31; A = 46045; D = 163.06;              E = 159.62;             F = 169.62 }

    @{ Row = 35; A = 46046; D = 157.83000000000001;  E = 149.26;             F = 158.26 }
    @{ Row = 36; A = 46045; D = 157.47999999999999;  E = 149.38999999999999; F = 158.38999999999999 }

    @{ Row = 40; A = 46046; D = 163.38;              E = 157.54;             F = 167.54 }
    @{ Row = 41; A = 46046; D = 163.1;               E = 157.97;             F = 167.97 }
    @{ Row = 42; A = 46045; D = 163.05000000000001;  E = 157.75;             F = 167.75 }
    @{ Row = 43; A = 46045; D = 162.76;              E = 158.16999999999999; F = 168.17 }

    @{ Row = 47; A = 46046; D = 157.33000000000001;  E = 151.07;             F = 161.07 }
    @{ Row = 48; A = 46046; D = 156.91;              E = 150.97;             F = 160.97 }
    @{ Row = 49; A = 46045; D = 157.1;               E = 151.18;             F = 161.18 }
    @{ Row = 50; A = 46045; D = 156.68;              E = 151.08000000000001; F = 161.08000000000001 }

    @{ Row = 54; A = 46046; D = 172.37;              E = 165.6;              F = 175.6 }
    @{ Row = 55; A = 46046; D = 165.17;              E = 163.46;             F = 173.46 }
    @{ Row = 56; A = 46046; D = 161.9 }
    @{ Row = 57; A = 46046; D = 162.24;              E = 157.88 }
    @{ Row = 58; A = 46046; D = 158.01;              E = 153.78;             F = 163.78 }
    @{ Row = 59; A = 46046; D = 164.66;              E = 163.6 }
    @{ Row = 60; A = 46045; D = 172.03;              E = 165.81;             F = 175.81 }
    @{ Row = 61; A = 46045; D = 164.86;              E = 163.58000000000001; F = 173.58 }
    @{ Row = 62; A = 46045; D = 161.55000000000001 }
    @{ Row = 63; A = 46045; D = 161.86000000000001;  E = 158 }
    @{ Row = 64; A = 46045; D = 157.63;              E = 153.9;              F = 163.9 }
    @{ Row = 65; A = 46045; D = 164.28;              E = 163.79 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("A")) { $ws.Range("A$r").Value = $u.A }
    if ($u.ContainsKey("D")) { $ws.Range("D$r").Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Range("F$r").Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Range("G$r").Value = $u.G }
}
